# Updated symbol list on Mon Dec 26 05:38:10 UTC 2022 with GitHub Actions
# Refresh the "Price" (column D) quotes and a couple of "Best in 24h" badge
# strings (column E) to match the latest crawl snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells are stored as text (not numbers) in the source sheet, so each
# new value is written with a leading apostrophe to force a text literal
# and keep the cell type consistent with the rest of the column.
$priceUpdates = @{
    "D2"  = "244.09"
    "D3"  = "22.99"
    "D4"  = "5.413"
    "D5"  = "0.05964"
    "D6"  = "3.460"
    "D7"  = "6.529"
    "D8"  = "0.8149"
    "D9"  = "0.9211"
    "D10" = "0.1412"
    "D11" = "0.07435"
    "D12" = "0.03265"
    "D13" = "0.03084"
    "D14" = "0.09356"
    "D15" = "3.845"
    "D16" = "0.001567"
    "D18" = "0.0005943"
    "D19" = "0.006081"
    "D20" = "0.004998"
    "D21" = "0.0009823"
    "D22" = "0.00007799"
    "D23" = "3.612"
    "D24" = "2.150"
    "D26" = "0.1324"
    "D27" = "0.0002396"
    "D40" = "0.03927"
    "D41" = "0.006172"
    "D42" = "0.1077"
    "D43" = "0.002620"
    "D44" = "0.006747"
    "D45" = "0.00005250"
    "D48" = "0.7805"
    "D50" = "0.00002101"
    "D51" = "0.0002001"
}

foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).Value = "'" + $priceUpdates[$addr]
}

# Column E "Best in 24h" badge text toggles for two rows.
$ws.Range("E20").Value = "19HotbitTokenHTBBestin24h"
$ws.Range("E44").Value = "43LocalTradersLCT"
